$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Floating-point precision refresh for existing timestamp cells (B8, B9)
$ws.Cells.Item(8, 2).Value = 46043.44030385416
$ws.Cells.Item(9, 2).Value = 46043.44086981482

# New row 10
$ws.Cells.Item(10, 1).Value = "lupa password master"
$ws.Cells.Item(10, 2).Value = 46045.62859663679
$ws.Cells.Item(10, 2).NumberFormat = $ws.Cells.Item(9, 2).NumberFormat
